$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates that don't involve row reordering.
# A leading apostrophe forces Excel to keep the value as text (matching the
# original inlineStr cell type) instead of auto-converting it to a number.
$ws.Range("D2").Value = "'243.38"
$ws.Range("D3").Value = "'23.18"
$ws.Range("D4").Value = "'5.403"
$ws.Range("D5").Value = "'0.05976"
$ws.Range("D6").Value = "'3.432"
$ws.Range("D7").Value = "'6.530"
$ws.Range("D8").Value = "'0.8088"
$ws.Range("D9").Value = "'0.9231"
$ws.Range("D11").Value = "'0.07414"
$ws.Range("D12").Value = "'0.03283"
$ws.Range("D13").Value = "'0.03091"
$ws.Range("D14").Value = "'0.09359"
$ws.Range("D15").Value = "'3.851"
$ws.Range("D16").Value = "'0.001573"
$ws.Range("D17").Value = "'0.04708"

# Row 18 (One / ONE) price + label update
$ws.Range("D18").Value = "'0.0005933"
$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").Value = "'0.005859"
$ws.Range("D20").Value = "'0.001278"
$ws.Range("D21").Value = "'0.004887"
$ws.Range("D22").Value = "'0.00006802"
$ws.Range("D24").Value = "'2.149"
$ws.Range("D27").Value = "'0.0002341"
$ws.Range("D40").Value = "'0.03972"

# Rows 41-43 were reordered: KickToken moved to row 41, BKEXToken to row 42,
# CEJI to row 43 (each with its own updated price/link/label).
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006392"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1077"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003201"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.009172"
$ws.Range("D45").Value = "'0.00005101"
$ws.Range("D47").Value = "'0.7004"
$ws.Range("D48").Value = "'0.002415"
$ws.Range("D49").Value = "'0.00002100"
